# Rename sheets and update the saved view (frozen pane top-left cell and
# active selection) on the "Base Classes" sheet (formerly "Sheet1").

$wb = $excel.ActiveWorkbook

# Rename worksheets
$wsBase = $wb.Worksheets.Item("Sheet1")
$wsBase.Name = "Base Classes"

$wsPrestige = $wb.Worksheets.Item("Sheet2")
$wsPrestige.Name = "Prestige Classes"

# Update the view on "Base Classes": scroll the frozen pane back to the top
# (B2) and move the active cell selection to B14.
$wsBase.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2

$wsBase.Range("B14").Select()
